$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.073.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.74%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.721.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.56%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'613.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +6.67%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'191.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +8.80%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.638"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.83%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  -0.35%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.719"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.43%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.35%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'58.00"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +9.03%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.0000290"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.43%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'10.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.28%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'4.312.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.88%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'3.718.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.49%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'19.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.37%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  -0.07%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  -0.35%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'12.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.65%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'68.853.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.64%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'411.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.44%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +1.39%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'89.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.14%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  -1.26%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'12.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.15%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'10.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.23%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  +1.05%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'3.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.11%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  +1.15%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'33.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.35%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'7.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -9.56%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'12.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.04%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'0.123"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +4.25%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'45.87"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").Value = "'630.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +4.80%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'65.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.69%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.415"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.99%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.0₃0822"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -11.85%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.13%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +0.17%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.141"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.22%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'3.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.82%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.0446"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.42%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  -0.50%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  +3.61%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'2.870.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +5.18%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  +1.30%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'9.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.33%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "'ApeXProtocol"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'3.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.11%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "'Monero"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'143.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.77%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'2.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -18.87%  "
$ws.Range("E51").Style = "Normal"
